$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = 123
$ws.Range("E1").Value = 456
$ws.Range("F1").Formula = "=D1+E1"
